$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.077.69'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '3.089.80'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.96'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '171.33'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.01%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '3.084.26'
$ws.Range('E8').Value = '  -1.52%  '
$ws.Range('E9').Value = '  -1.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.45'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.96%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.152'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.73%  '
$ws.Range('E12').Value = '  -1.78%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000245'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.45'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.73%  '
$ws.Range('E15').Value = '  -1.75%  '
$ws.Range('D16').Value = '3.607.06'
$ws.Range('E16').Value = '  -1.36%  '
$ws.Range('D17').Value = '67.019.21'
$ws.Range('E17').Value = '  -0.37%  '
$ws.Range('E18').Value = '  -1.79%  '
$ws.Range('D19').Value = '3.090.78'
$ws.Range('E19').Value = '  -1.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.50'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '487.15'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.694'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.75'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.82%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.48'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.93'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.58%  '
$ws.Range('E26').Value = '  -3.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.41'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.34%  '
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.78'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.58%  '
$ws.Range('E30').Value = '  -3.93%  '
$ws.Range('E31').Value = '  -1.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.18'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.34%  '
$ws.Range('E33').Value = '  -2.10%  '
$ws.Range('D34').Value = '0.0₃0941'
$ws.Range('E34').Value = '  -5.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('E36').Value = '  -3.28%  '
$ws.Range('E37').Value = '  -2.80%  '
$ws.Range('E38').Value = '  -2.82%  '
$ws.Range('E39').Value = '  -4.90%  '
$ws.Range('E40').Value = '  +0.17%  '
$ws.Range('E41').Value = '  -2.94%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.40'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Value = '2.786.38'
$ws.Range('E43').Value = '  -2.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '379.32'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.27%  '
$ws.Range('E45').Value = '  -8.28%  '
$ws.Range('E46').Value = '  -2.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '134.95'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '24.74'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.80%  '
$ws.Range('E50').Value = '  -1.93%  '
$ws.Range('E51').Value = '  -2.24%  '
